$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A22").Value = "Find if there is a sub array with 0 sum"
$ws.Range("B22").Value = "SubArrayWIth0Sum"

$ws.Range("B22").Select()
